$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.197.97"
$ws.Range("E2").Value = "  +2.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.321.99"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "544.69"
$ws.Range("E5").Value = "  +1.02%  "

# Row 6 - Solana
$ws.Range("D6").Value = "130.99"
$ws.Range("E6").Value = "  -1.31%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  -1.38%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.320.20"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +0.35%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "5.52"
$ws.Range("E11").Value = "  -0.49%  "

# Row 12 - TRON (unchanged)

# Row 13 - Cardano
$ws.Range("D13").Value = "0.333"
$ws.Range("E13").Value = "  -0.48%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "23.54"
$ws.Range("E14").Value = "  -1.29%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.739.80"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "60.155.24"
$ws.Range("E16").Value = "  +2.52%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -0.70%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.322.45"
$ws.Range("E18").Value = "  +1.95%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "10.57"
$ws.Range("E19").Value = "  -0.53%  "

# Row 20 - Polkadot
$ws.Range("D20").Value = "4.12"
$ws.Range("E20").Value = "  -2.19%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "314.18"
$ws.Range("E21").Value = "  -0.64%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.62"
$ws.Range("E22").Value = "  -0.32%  "

# Row 23 - Dai
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "63.89"
$ws.Range("E24").Value = "  +1.07%  "

# Row 25 - Kaspa
$ws.Range("D25").Value = "0.171"
$ws.Range("E25").Value = "  -0.95%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.24%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "7.80"
$ws.Range("E27").Value = "  -1.88%  "

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  +4.98%  "

# Row 29 - SuiNetwork
$ws.Range("E29").Value = "  +6.19%  "

# Row 30 - Monero
$ws.Range("D30").Value = "172.51"
$ws.Range("E30").Value = "  +0.67%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.31%  "

# Row 32 - PEPE
$ws.Range("D32").Value = "0.0₃0730"
$ws.Range("E32").Value = "  -0.69%  "

# Row 33 - Aptos
$ws.Range("D33").Value = "5.93"
$ws.Range("E33").Value = "  +0.94%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  +8.26%  "

# Row 35 - PolygonEcosystemToken
$ws.Range("D35").Value = "0.379"
$ws.Range("E35").Value = "  -1.95%  "

# Row 36 - USDe
$ws.Range("E36").Value = "  +0.00%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "17.78"
$ws.Range("E37").Value = "  -1.13%  "

# Row 38 - FirstDigitalUSD
$ws.Range("E38").Value = "  -0.03%  "

# Row 39 - NEARProtocol
$ws.Range("D39").Value = "4.04"
$ws.Range("E39").Value = "  +0.27%  "

# Row 40 - Bittensor
$ws.Range("D40").Value = "318.42"
$ws.Range("E40").Value = "  +7.04%  "

# Row 41 - was Stacks, now OKB
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "37.86"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42 - was OKB, now Stacks
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.52"
$ws.Range("E42").Value = "  +0.02%  "

# Row 43 - Aave
$ws.Range("D43").Value = "137.37"
$ws.Range("E43").Value = "  -2.78%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "3.47"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45 - Stellar
$ws.Range("D45").Value = "0.0938"
$ws.Range("E45").Value = "  -1.85%  "

# Row 46 - InjectiveProtocol
$ws.Range("D46").Value = "19.01"
$ws.Range("E46").Value = "  +3.00%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "0.562"
$ws.Range("E47").Value = "  +0.87%  "

# Row 48 - Hedera
$ws.Range("D48").Value = "0.0493"
$ws.Range("E48").Value = "  -0.65%  "

# Row 49 - was VeChain, now BabyDogeCoin
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0220"
$ws.Range("E49").Value = "  +17.50%  "

# Row 50 - was BabyDogeCoin, now VeChain
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0213"
$ws.Range("E50").Value = "  +0.69%  "

# Row 51 - WhiteBITCoin
$ws.Range("D51").Value = "11.01"
$ws.Range("E51").Value = "  +0.24%  "
